# Update column G ("K" = strikeouts) values for rows 2-8 on Sheet1.
# These values were regenerated as part of switching save_data generation
# to use K instead of Strike#.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 1
